$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark from the "Leer archivo de disco duro (JSON)." paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Merge the two runs of the "Actualizar" bullet into a single run.
$d.Content.Find.Execute("Actualizar (archivo o carpeta sobre el grafo).", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "Actualizar (archivo o carpeta sobre el grafo).", 2) | Out-Null

# 3. Merge the two runs of the "Eliminar" bullet into a single run.
$d.Content.Find.Execute("Eliminar (archivo o carpeta sobre el grafo).", $true, $false, $false,
                         $false, $false, $true, 1, $false,
                         "Eliminar (archivo o carpeta sobre el grafo).", 2) | Out-Null

# 4. Add a new bullet "Balanceo del grafo." at the end of the list, carrying the
#    "_GoBack" bookmark that used to sit on the first bullet.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)

# A placeholder trailing character is appended (and removed afterwards) so the
# bookmark's insertion point is never the very last character position of the
# paragraph -- inserting a zero-length bookmark exactly there mis-locates it.
$newPara.Range.InsertAfter("Balanceo del grafo.X")

$bmPos = $newPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()
